# Append 20 new daily NAV rows (2024-09-02 .. 2024-09-27) after the
# existing last row (654), extending the used range to A1:J674.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 655 (2024-09-02)
$ws.Cells.Item(655, 1).NumberFormat = "@"
$ws.Cells.Item(655, 1).Value = "2024-09-02"
$ws.Cells.Item(655, 1).ClearFormats()
$ws.Cells.Item(655, 3).Value = 1840.550048828125
$ws.Cells.Item(655, 4).Value = 608.5800170898438
$ws.Cells.Item(655, 5).Value = 1111.550048828125
$ws.Cells.Item(655, 6).Value = 177.5399932861328
$ws.Cells.Item(655, 7).Value = 670.2000122070312
$ws.Cells.Item(655, 8).Value = 22669.95040893555
$ws.Cells.Item(655, 9).Value = 0
$ws.Cells.Item(655, 10).Value = 220.09064177077

# Row 656 (2024-09-03)
$ws.Cells.Item(656, 1).NumberFormat = "@"
$ws.Cells.Item(656, 1).Value = "2024-09-03"
$ws.Cells.Item(656, 1).ClearFormats()
$ws.Cells.Item(656, 3).Value = 1865.599975585938
$ws.Cells.Item(656, 4).Value = 599.9400024414062
$ws.Cells.Item(656, 5).Value = 1114
$ws.Cells.Item(656, 6).Value = 178.4600067138672
$ws.Cells.Item(656, 7).Value = 659.0999755859375
$ws.Cells.Item(656, 8).Value = 22746.27993774414
$ws.Cells.Item(656, 9).Value = 0.003366991432787071
$ws.Cells.Item(656, 10).Value = 220.8316850760488

# Row 657 (2024-09-04)
$ws.Cells.Item(657, 1).NumberFormat = "@"
$ws.Cells.Item(657, 1).Value = "2024-09-04"
$ws.Cells.Item(657, 1).ClearFormats()
$ws.Cells.Item(657, 3).Value = 1871.900024414062
$ws.Cells.Item(657, 4).Value = 609
$ws.Cells.Item(657, 5).Value = 1127.900024414062
$ws.Cells.Item(657, 6).Value = 176.0200042724609
$ws.Cells.Item(657, 7).Value = 650.8499755859375
$ws.Cells.Item(657, 8).Value = 22782.94021606445
$ws.Cells.Item(657, 9).Value = 0.001611704349926693
$ws.Cells.Item(657, 10).Value = 221.1876004634875

# Row 658 (2024-09-05)
$ws.Cells.Item(658, 1).NumberFormat = "@"
$ws.Cells.Item(658, 1).Value = "2024-09-05"
$ws.Cells.Item(658, 1).ClearFormats()
$ws.Cells.Item(658, 3).Value = 1864.949951171875
$ws.Cells.Item(658, 4).Value = 602.1799926757812
$ws.Cells.Item(658, 5).Value = 1115.150024414062
$ws.Cells.Item(658, 6).Value = 173.4799957275391
$ws.Cells.Item(658, 7).Value = 643.8499755859375
$ws.Cells.Item(658, 8).Value = 22586.02963256836
$ws.Cells.Item(658, 9).Value = -0.008642896027846764
$ws.Cells.Item(658, 10).Value = 219.2758990300327

# Row 659 (2024-09-06)
$ws.Cells.Item(659, 1).NumberFormat = "@"
$ws.Cells.Item(659, 1).Value = "2024-09-06"
$ws.Cells.Item(659, 1).ClearFormats()
$ws.Cells.Item(659, 3).Value = 1857.150024414062
$ws.Cells.Item(659, 4).Value = 597.2999877929688
$ws.Cells.Item(659, 5).Value = 1100
$ws.Cells.Item(659, 6).Value = 169.8500061035156
$ws.Cells.Item(659, 7).Value = 634.7000122070312
$ws.Cells.Item(659, 8).Value = 22350.45025634766
$ws.Cells.Item(659, 9).Value = -0.01043031378480992
$ws.Cells.Item(659, 10).Value = 216.9887825977031

# Row 660 (2024-09-09)
$ws.Cells.Item(660, 1).NumberFormat = "@"
$ws.Cells.Item(660, 1).Value = "2024-09-09"
$ws.Cells.Item(660, 1).ClearFormats()
$ws.Cells.Item(660, 3).Value = 1860.449951171875
$ws.Cells.Item(660, 4).Value = 610.3400268554688
$ws.Cells.Item(660, 5).Value = 1104.150024414062
$ws.Cells.Item(660, 6).Value = 168.3300018310547
$ws.Cells.Item(660, 7).Value = 635.2000122070312
$ws.Cells.Item(660, 8).Value = 22404.27005004883
$ws.Cells.Item(660, 9).Value = 0.002407995950143633
$ws.Cells.Item(660, 10).Value = 217.511290707425

# Row 661 (2024-09-10)
$ws.Cells.Item(661, 1).NumberFormat = "@"
$ws.Cells.Item(661, 1).Value = "2024-09-10"
$ws.Cells.Item(661, 1).ClearFormats()
$ws.Cells.Item(661, 3).Value = 1824.5
$ws.Cells.Item(661, 4).Value = 608
$ws.Cells.Item(661, 5).Value = 1113.199951171875
$ws.Cells.Item(661, 6).Value = 169.75
$ws.Cells.Item(661, 7).Value = 637.0499877929688
$ws.Cells.Item(661, 8).Value = 22289.99975585938
$ws.Cells.Item(661, 9).Value = -0.005100380147810443
$ws.Cells.Item(661, 10).Value = 216.4019004383762

# Row 662 (2024-09-11)
$ws.Cells.Item(662, 1).NumberFormat = "@"
$ws.Cells.Item(662, 1).Value = "2024-09-11"
$ws.Cells.Item(662, 1).ClearFormats()
$ws.Cells.Item(662, 3).Value = 1833.150024414062
$ws.Cells.Item(662, 4).Value = 627.6599731445312
$ws.Cells.Item(662, 5).Value = 1112.599975585938
$ws.Cells.Item(662, 6).Value = 165.8800048828125
$ws.Cells.Item(662, 7).Value = 627.2000122070312
$ws.Cells.Item(662, 8).Value = 22284.95007324219
$ws.Cells.Item(662, 9).Value = -0.0002265447587481507
$ws.Cells.Item(662, 10).Value = 216.3528757220488

# Row 663 (2024-09-12)
$ws.Cells.Item(663, 1).NumberFormat = "@"
$ws.Cells.Item(663, 1).Value = "2024-09-12"
$ws.Cells.Item(663, 1).ClearFormats()
$ws.Cells.Item(663, 3).Value = 1854.849975585938
$ws.Cells.Item(663, 4).Value = 645.5999755859375
$ws.Cells.Item(663, 5).Value = 1120.099975585938
$ws.Cells.Item(663, 6).Value = 167.0200042724609
$ws.Cells.Item(663, 7).Value = 651.0999755859375
$ws.Cells.Item(663, 8).Value = 22615.88967895508
$ws.Cells.Item(663, 9).Value = 0.01485036334500268
$ws.Cells.Item(663, 10).Value = 219.5657945372574

# Row 664 (2024-09-13)
$ws.Cells.Item(664, 1).NumberFormat = "@"
$ws.Cells.Item(664, 1).Value = "2024-09-13"
$ws.Cells.Item(664, 1).ClearFormats()
$ws.Cells.Item(664, 3).Value = 1894.449951171875
$ws.Cells.Item(664, 4).Value = 646.6500244140625
$ws.Cells.Item(664, 5).Value = 1118.550048828125
$ws.Cells.Item(664, 6).Value = 167.25
$ws.Cells.Item(664, 7).Value = 633.4500122070312
$ws.Cells.Item(664, 8).Value = 22746.35009765625
$ws.Cells.Item(664, 9).Value = 0.005768529142701387
$ws.Cells.Item(664, 10).Value = 220.832366221786

# Row 665 (2024-09-16)
$ws.Cells.Item(665, 1).NumberFormat = "@"
$ws.Cells.Item(665, 1).Value = "2024-09-16"
$ws.Cells.Item(665, 1).ClearFormats()
$ws.Cells.Item(665, 3).Value = 1857.599975585938
$ws.Cells.Item(665, 4).Value = 621.0499877929688
$ws.Cells.Item(665, 5).Value = 1115.849975585938
$ws.Cells.Item(665, 6).Value = 163.9600067138672
$ws.Cells.Item(665, 7).Value = 665.9500122070312
$ws.Cells.Item(665, 8).Value = 22506.51992797852
$ws.Cells.Item(665, 9).Value = -0.01054367705799297
$ws.Cells.Item(665, 10).Value = 218.5039810683911

# Row 666 (2024-09-17)
$ws.Cells.Item(666, 1).NumberFormat = "@"
$ws.Cells.Item(666, 1).Value = "2024-09-17"
$ws.Cells.Item(666, 1).ClearFormats()
$ws.Cells.Item(666, 3).Value = 1848.699951171875
$ws.Cells.Item(666, 4).Value = 649.6500244140625
$ws.Cells.Item(666, 5).Value = 1110.949951171875
$ws.Cells.Item(666, 6).Value = 160.6000061035156
$ws.Cells.Item(666, 7).Value = 666.3499755859375
$ws.Cells.Item(666, 8).Value = 22484.49969482422
$ws.Cells.Item(666, 9).Value = -0.0009783935155129372
$ws.Cells.Item(666, 10).Value = 218.2901981902

# Row 667 (2024-09-18)
$ws.Cells.Item(667, 1).NumberFormat = "@"
$ws.Cells.Item(667, 1).Value = "2024-09-18"
$ws.Cells.Item(667, 1).ClearFormats()
$ws.Cells.Item(667, 3).Value = 1888.199951171875
$ws.Cells.Item(667, 4).Value = 646.7000122070312
$ws.Cells.Item(667, 5).Value = 1079.949951171875
$ws.Cells.Item(667, 6).Value = 158.5599975585938
$ws.Cells.Item(667, 7).Value = 651.7000122070312
$ws.Cells.Item(667, 8).Value = 22442.71960449219
$ws.Cells.Item(667, 9).Value = -0.00185817300358472
$ws.Cells.Item(667, 10).Value = 217.8845772369758

# Row 668 (2024-09-19)
$ws.Cells.Item(668, 1).NumberFormat = "@"
$ws.Cells.Item(668, 1).Value = "2024-09-19"
$ws.Cells.Item(668, 1).ClearFormats()
$ws.Cells.Item(668, 3).Value = 1890.400024414062
$ws.Cells.Item(668, 4).Value = 652.1500244140625
$ws.Cells.Item(668, 5).Value = 1054.449951171875
$ws.Cells.Item(668, 6).Value = 155.25
$ws.Cells.Item(668, 7).Value = 649.5999755859375
$ws.Cells.Item(668, 8).Value = 22292.29992675781
$ws.Cells.Item(668, 9).Value = -0.006702381903139165
$ws.Cells.Item(668, 10).Value = 216.4242315895295

# Row 669 (2024-09-20)
$ws.Cells.Item(669, 1).NumberFormat = "@"
$ws.Cells.Item(669, 1).Value = "2024-09-20"
$ws.Cells.Item(669, 1).ClearFormats()
$ws.Cells.Item(669, 3).Value = 1916.800048828125
$ws.Cells.Item(669, 4).Value = 654.4500122070312
$ws.Cells.Item(669, 5).Value = 1054.599975585938
$ws.Cells.Item(669, 6).Value = 161.4299926757812
$ws.Cells.Item(669, 7).Value = 665.1500244140625
$ws.Cells.Item(669, 8).Value = 22632.26013183594
$ws.Cells.Item(669, 9).Value = 0.01525011803156592
$ws.Cells.Item(669, 10).Value = 219.7247266661608

# Row 670 (2024-09-23)
$ws.Cells.Item(670, 1).NumberFormat = "@"
$ws.Cells.Item(670, 1).Value = "2024-09-23"
$ws.Cells.Item(670, 1).ClearFormats()
$ws.Cells.Item(670, 3).Value = 1919.949951171875
$ws.Cells.Item(670, 4).Value = 654.0999755859375
$ws.Cells.Item(670, 5).Value = 1055.25
$ws.Cells.Item(670, 6).Value = 159.5599975585938
$ws.Cells.Item(670, 7).Value = 672
$ws.Cells.Item(670, 8).Value = 22635.46960449219
$ws.Cells.Item(670, 9).Value = 0.0001418096397599883
$ws.Cells.Item(670, 10).Value = 219.7558857504957

# Row 671 (2024-09-24)
$ws.Cells.Item(671, 1).NumberFormat = "@"
$ws.Cells.Item(671, 1).Value = "2024-09-24"
$ws.Cells.Item(671, 1).ClearFormats()
$ws.Cells.Item(671, 3).Value = 1904.650024414062
$ws.Cells.Item(671, 4).Value = 646.8499755859375
$ws.Cells.Item(671, 5).Value = 1051.550048828125
$ws.Cells.Item(671, 6).Value = 158.7400054931641
$ws.Cells.Item(671, 7).Value = 675.25
$ws.Cells.Item(671, 8).Value = 22510.13034057617
$ws.Cells.Item(671, 9).Value = -0.005537294613544976
$ws.Cells.Item(671, 10).Value = 218.5390326680347

# Row 672 (2024-09-25)
$ws.Cells.Item(672, 1).NumberFormat = "@"
$ws.Cells.Item(672, 1).Value = "2024-09-25"
$ws.Cells.Item(672, 1).ClearFormats()
$ws.Cells.Item(672, 3).Value = 1928.5
$ws.Cells.Item(672, 4).Value = 633.2999877929688
$ws.Cells.Item(672, 5).Value = 1063.449951171875
$ws.Cells.Item(672, 6).Value = 156.9400024414062
$ws.Cells.Item(672, 7).Value = 667.3499755859375
$ws.Cells.Item(672, 8).Value = 22551.57971191406
$ws.Cells.Item(672, 9).Value = 0.001841365230265907
$ws.Cells.Item(672, 10).Value = 218.9414428442456

# Row 673 (2024-09-26)
$ws.Cells.Item(673, 1).NumberFormat = "@"
$ws.Cells.Item(673, 1).Value = "2024-09-26"
$ws.Cells.Item(673, 1).ClearFormats()
$ws.Cells.Item(673, 3).Value = 1982.800048828125
$ws.Cells.Item(673, 4).Value = 626.8499755859375
$ws.Cells.Item(673, 5).Value = 1068
$ws.Cells.Item(673, 6).Value = 156.8500061035156
$ws.Cells.Item(673, 7).Value = 665.3499755859375
$ws.Cells.Item(673, 8).Value = 22805.50018310547
$ws.Cells.Item(673, 9).Value = 0.01125954254358773
$ws.Cells.Item(673, 10).Value = 221.4066233345048

# Row 674 (2024-09-27)
$ws.Cells.Item(674, 1).NumberFormat = "@"
$ws.Cells.Item(674, 1).Value = "2024-09-27"
$ws.Cells.Item(674, 1).ClearFormats()
$ws.Cells.Item(674, 3).Value = 2010.699951171875
$ws.Cells.Item(674, 4).Value = 608.5499877929688
$ws.Cells.Item(674, 5).Value = 1075.949951171875
$ws.Cells.Item(674, 6).Value = 156.8099975585938
$ws.Cells.Item(674, 7).Value = 654.2999877929688
$ws.Cells.Item(674, 8).Value = 22858.51940917969
$ws.Cells.Item(674, 9).Value = 0.002324843816120109
$ws.Cells.Item(674, 10).Value = 221.9213591536121
